$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Planning & design" update: record the actual hours burned each week.
# These sit on the diagonal of the weekly-tracking grid (row = sprint,
# column = week); everything else in that grid stays at 0.
$ws.Range("C2").Value = 4      # Sprint 1 / Week 1
$ws.Range("D3").Value = 4.5    # Sprint 2 / Week 2
$ws.Range("E4").Value = 5.5    # Sprint 3 / Week 3
$ws.Range("F5").Value = 6      # Sprint 4 / Week 4
$ws.Range("G6").Value = 3      # Sprint 5 / Week 5

# Row 20 (Actual Hours) and row 21 (Remaining Effort) are SUM/IF formulas
# over the grid above, so they - and the chart series that read them -
# recalculate automatically once the inputs change.

# Match the author's final viewport/selection state on save.
$ws.Activate()
$ws.Range("G35").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
